# CIERRE 26 DIC 23
# Advance the payroll receipt workbook from "SEMANA 49 (DEL 04 AL 10 DICIEMBRE 2023)"
# to "SEMANA 50 (DEL 11 AL 17 DICIEMBRE 2023)" and update the period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the week/period label (shared by several cells through formulas
# referencing B9, e.g. H9=B9, B28=B9, H28=B28, B46=H28, B62=B46).
$ws.Range("B9").Value = "SEMANA  50       DEL    11     Al    17   DICIEMBRE     2023"

# --- Left table (employee #1) ---
$ws.Range("K4").Value = 0

# --- Right table (employee #2), row 23 block ---
$ws.Range("J23").Value = 5
$ws.Range("K23").Value = 2167
$ws.Range("E25").Value = 0
$ws.Range("K25").Value = 500

# Reset the view: scroll back to the top and select B11 (matches the
# saved sheetView - no topLeftCell override, activeCell/sqref = B11).
$ws.Range("B11").Select()
